$wb = $excel.ActiveWorkbook

# --- Sheet 1: append new portfolio snapshot rows ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("A51:C60").NumberFormat = "@"
$ws.Range("I51:I60").NumberFormat = "@"

$ws.Range("A51").Value = '大智 (稳健智远)'
$ws.Range("B51").Value = '000333'
$ws.Range("C51").Value = '美的集团'
$ws.Range("D51").Value = 2.77
$ws.Range("E51").Value = 39.71497483560155
$ws.Range("F51").Value = 71.48
$ws.Range("G51").Value = 2838.826401248799
$ws.Range("H51").Value = 102504.055170665
$ws.Range("I51").Value = '202506131600'

$ws.Range("A52").Value = '大智 (稳健智远)'
$ws.Range("B52").Value = '510050'
$ws.Range("C52").Value = '上证50ETF'
$ws.Range("D52").Value = 4.85
$ws.Range("E52").Value = 1808.278443601665
$ws.Range("F52").Value = 2.75
$ws.Range("G52").Value = 4972.765719904579
$ws.Range("H52").Value = 102504.055170665
$ws.Range("I52").Value = '202506131600'

$ws.Range("A53").Value = '大智 (稳健智远)'
$ws.Range("B53").Value = '510300'
$ws.Range("C53").Value = '沪深300ETF'
$ws.Range("D53").Value = 4.84
$ws.Range("E53").Value = 1247.712126085149
$ws.Range("F53").Value = 3.98
$ws.Range("G53").Value = 4965.894261818893
$ws.Range("H53").Value = 102504.055170665
$ws.Range("I53").Value = '202506131600'

$ws.Range("A54").Value = '大智 (稳健智远)'
$ws.Range("B54").Value = '518880'
$ws.Range("C54").Value = '黄金ETF'
$ws.Range("D54").Value = 4.99
$ws.Range("E54").Value = 674.4389870730533
$ws.Range("F54").Value = 7.59
$ws.Range("G54").Value = 5118.991911884475
$ws.Range("H54").Value = 102504.055170665
$ws.Range("I54").Value = '202506131600'

$ws.Range("A55").Value = '大智 (稳健智远)'
$ws.Range("B55").Value = '600085'
$ws.Range("C55").Value = '同仁堂'
$ws.Range("D55").Value = 1.9
$ws.Range("E55").Value = 52.96734947562319
$ws.Range("F55").Value = 36.7
$ws.Range("G55").Value = 1943.901725755371
$ws.Range("H55").Value = 102504.055170665
$ws.Range("I55").Value = '202506131600'

$ws.Range("A56").Value = '大智 (稳健智远)'
$ws.Range("B56").Value = '600900'
$ws.Range("C56").Value = '长江电力'
$ws.Range("D56").Value = 19.89
$ws.Range("E56").Value = 665.8903941748626
$ws.Range("F56").Value = 30.62
$ws.Range("G56").Value = 20389.56386963429
$ws.Range("H56").Value = 102504.055170665
$ws.Range("I56").Value = '202506131600'

$ws.Range("A57").Value = '大智 (稳健智远)'
$ws.Range("B57").Value = '600989'
$ws.Range("C57").Value = '宝丰能源'
$ws.Range("D57").Value = 4.97
$ws.Range("E57").Value = 306.7515982999751
$ws.Range("F57").Value = 16.61
$ws.Range("G57").Value = 5095.144047762587
$ws.Range("H57").Value = 102504.055170665
$ws.Range("I57").Value = '202506131600'

$ws.Range("A58").Value = '大智 (稳健智远)'
$ws.Range("B58").Value = 'HK02899'
$ws.Range("C58").Value = '紫金矿业'
$ws.Range("D58").Value = 21.7
$ws.Range("E58").Value = 1106.618293645365
$ws.Range("F58").Value = 20.1
$ws.Range("G58").Value = 22243.02770227184
$ws.Range("H58").Value = 102504.055170665
$ws.Range("I58").Value = '202506131600'

$ws.Range("A59").Value = '大智 (稳健智远)'
$ws.Range("B59").Value = 'HK06881'
$ws.Range("C59").Value = '中国银河'
$ws.Range("D59").Value = 4.87
$ws.Range("E59").Value = 610.1281790147427
$ws.Range("F59").Value = 8.18
$ws.Range("G59").Value = 4990.848504340595
$ws.Range("H59").Value = 102504.055170665
$ws.Range("I59").Value = '202506131600'

$ws.Range("A60").Value = '大智 (稳健智远)'
$ws.Range("B60").Value = '100000'
$ws.Range("C60").Value = '现金'
$ws.Range("D60").Value = 29.21
$ws.Range("E60").Value = 29945.09102604357
$ws.Range("F60").Value = 1
$ws.Range("G60").Value = 29945.09102604357
$ws.Range("H60").Value = 102504.055170665
$ws.Range("I60").Value = '202506131600'

# --- Sheet 2: append new portfolio snapshot rows ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("A37:C43").NumberFormat = "@"
$ws.Range("I37:I43").NumberFormat = "@"

$ws.Range("A37").Value = '大成 (锐进先锋)'
$ws.Range("B37").Value = '000725'
$ws.Range("C37").Value = '京东方A'
$ws.Range("D37").Value = 4.91
$ws.Range("E37").Value = 1243.91705951017
$ws.Range("F37").Value = 3.89
$ws.Range("G37").Value = 4838.837361494561
$ws.Range("H37").Value = 98482.59985089369
$ws.Range("I37").Value = '202506131600'

$ws.Range("A38").Value = '大成 (锐进先锋)'
$ws.Range("B38").Value = '159781'
$ws.Range("C38").Value = '科创创业ETF'
$ws.Range("D38").Value = 4.91
$ws.Range("E38").Value = 9122.058436407913
$ws.Range("F38").Value = 0.53
$ws.Range("G38").Value = 4834.690971296194
$ws.Range("H38").Value = 98482.59985089369
$ws.Range("I38").Value = '202506131600'

$ws.Range("A39").Value = '大成 (锐进先锋)'
$ws.Range("B39").Value = '513100'
$ws.Range("C39").Value = '纳指ETF'
$ws.Range("D39").Value = 4.94
$ws.Range("E39").Value = 3137.523283860047
$ws.Range("F39").Value = 1.55
$ws.Range("G39").Value = 4863.161089983073
$ws.Range("H39").Value = 98482.59985089369
$ws.Range("I39").Value = '202506131600'

$ws.Range("A40").Value = '大成 (锐进先锋)'
$ws.Range("B40").Value = '513290'
$ws.Range("C40").Value = '纳指生物科技ETF'
$ws.Range("D40").Value = 0.99
$ws.Range("E40").Value = 879.627063510763
$ws.Range("F40").Value = 1.11
$ws.Range("G40").Value = 976.3860404969471
$ws.Range("H40").Value = 98482.59985089369
$ws.Range("I40").Value = '202506131600'

$ws.Range("A41").Value = '大成 (锐进先锋)'
$ws.Range("B41").Value = '603119'
$ws.Range("C41").Value = '浙江荣泰'
$ws.Range("D41").Value = 44.66
$ws.Range("E41").Value = 1051.546584462582
$ws.Range("F41").Value = 41.83
$ws.Range("G41").Value = 43986.1936280698
$ws.Range("H41").Value = 98482.59985089369
$ws.Range("I41").Value = '202506131600'

$ws.Range("A42").Value = '大成 (锐进先锋)'
$ws.Range("B42").Value = '688290'
$ws.Range("C42").Value = '景业智能'
$ws.Range("D42").Value = 9.57
$ws.Range("E42").Value = 161.2702946560293
$ws.Range("F42").Value = 58.46
$ws.Range("G42").Value = 9427.861425591473
$ws.Range("H42").Value = 98482.59985089369
$ws.Range("I42").Value = '202506131600'

$ws.Range("A43").Value = '大成 (锐进先锋)'
$ws.Range("B43").Value = '100000'
$ws.Range("C43").Value = '现金'
$ws.Range("D43").Value = 30.01
$ws.Range("E43").Value = 29555.46933396164
$ws.Range("F43").Value = 1
$ws.Range("G43").Value = 29555.46933396164
$ws.Range("H43").Value = 98482.59985089369
$ws.Range("I43").Value = '202506131600'

# --- Sheet 3: append new portfolio snapshot rows ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("A76:C90").NumberFormat = "@"
$ws.Range("I76:I90").NumberFormat = "@"

$ws.Range("A76").Value = '范式进化投资组合'
$ws.Range("B76").Value = '000333'
$ws.Range("C76").Value = '美的集团'
$ws.Range("D76").Value = 0.95
$ws.Range("E76").Value = 13.2658076636599
$ws.Range("F76").Value = 71.48
$ws.Range("G76").Value = 948.2399317984097
$ws.Range("H76").Value = 99958.03107115487
$ws.Range("I76").Value = '202506131600'

$ws.Range("A77").Value = '范式进化投资组合'
$ws.Range("B77").Value = '000725'
$ws.Range("C77").Value = '京东方A'
$ws.Range("D77").Value = 4.91
$ws.Range("E77").Value = 1262.931689191865
$ws.Range("F77").Value = 3.89
$ws.Range("G77").Value = 4912.804270956355
$ws.Range("H77").Value = 99958.03107115487
$ws.Range("I77").Value = '202506131600'

$ws.Range("A78").Value = '范式进化投资组合'
$ws.Range("B78").Value = '159781'
$ws.Range("C78").Value = '科创创业ETF'
$ws.Range("D78").Value = 4.91
$ws.Range("E78").Value = 9261.499054073673
$ws.Range("F78").Value = 0.53
$ws.Range("G78").Value = 4908.594498659047
$ws.Range("H78").Value = 99958.03107115487
$ws.Range("I78").Value = '202506131600'

$ws.Range("A79").Value = '范式进化投资组合'
$ws.Range("B79").Value = '510050'
$ws.Range("C79").Value = '上证50ETF'
$ws.Range("D79").Value = 4.99
$ws.Range("E79").Value = 1812.03242362311
$ws.Range("F79").Value = 2.75
$ws.Range("G79").Value = 4983.089164963552
$ws.Range("H79").Value = 99958.03107115487
$ws.Range("I79").Value = '202506131600'

$ws.Range("A80").Value = '范式进化投资组合'
$ws.Range("B80").Value = '510300'
$ws.Range("C80").Value = '沪深300ETF'
$ws.Range("D80").Value = 4.98
$ws.Range("E80").Value = 1250.302372299946
$ws.Range("F80").Value = 3.98
$ws.Range("G80").Value = 4976.203441753785
$ws.Range("H80").Value = 99958.03107115487
$ws.Range("I80").Value = '202506131600'

$ws.Range("A81").Value = '范式进化投资组合'
$ws.Range("B81").Value = '513100'
$ws.Range("C81").Value = '纳指ETF'
$ws.Range("D81").Value = 0.99
$ws.Range("E81").Value = 637.096750216533
$ws.Range("F81").Value = 1.55
$ws.Range("G81").Value = 987.4999628356262
$ws.Range("H81").Value = 99958.03107115487
$ws.Range("I81").Value = '202506131600'

$ws.Range("A82").Value = '范式进化投资组合'
$ws.Range("B82").Value = '513290'
$ws.Range("C82").Value = '纳指生物科技ETF'
$ws.Range("D82").Value = 0.99
$ws.Range("E82").Value = 893.0731230713899
$ws.Range("F82").Value = 1.11
$ws.Range("G82").Value = 991.3111666092428
$ws.Range("H82").Value = 99958.03107115487
$ws.Range("I82").Value = '202506131600'

$ws.Range("A83").Value = '范式进化投资组合'
$ws.Range("B83").Value = '518880'
$ws.Range("C83").Value = '黄金ETF'
$ws.Range("D83").Value = 1.03
$ws.Range("E83").Value = 135.1678240324266
$ws.Range("F83").Value = 7.59
$ws.Range("G83").Value = 1025.923784406118
$ws.Range("H83").Value = 99958.03107115487
$ws.Range("I83").Value = '202506131600'

$ws.Range("A84").Value = '范式进化投资组合'
$ws.Range("B84").Value = '600085'
$ws.Range("C84").Value = '同仁堂'
$ws.Range("D84").Value = 0.97
$ws.Range("E84").Value = 26.53865475829018
$ws.Range("F84").Value = 36.7
$ws.Range("G84").Value = 973.9686296292497
$ws.Range("H84").Value = 99958.03107115487
$ws.Range("I84").Value = '202506131600'

$ws.Range("A85").Value = '范式进化投资组合'
$ws.Range("B85").Value = '600900'
$ws.Range("C85").Value = '长江电力'
$ws.Range("D85").Value = 1.02
$ws.Range("E85").Value = 33.36363902067901
$ws.Range("F85").Value = 30.62
$ws.Range("G85").Value = 1021.594626813191
$ws.Range("H85").Value = 99958.03107115487
$ws.Range("I85").Value = '202506131600'

$ws.Range("A86").Value = '范式进化投资组合'
$ws.Range("B86").Value = '600989'
$ws.Range("C86").Value = '宝丰能源'
$ws.Range("D86").Value = 5.11
$ws.Range("E86").Value = 307.3884135955614
$ws.Range("F86").Value = 16.61
$ws.Range("G86").Value = 5105.721549822275
$ws.Range("H86").Value = 99958.03107115487
$ws.Range("I86").Value = '202506131600'

$ws.Range("A87").Value = '范式进化投资组合'
$ws.Range("B87").Value = '603119'
$ws.Range("C87").Value = '浙江荣泰'
$ws.Range("D87").Value = 0.99
$ws.Range("E87").Value = 23.72490270018873
$ws.Range("F87").Value = 41.83
$ws.Range("G87").Value = 992.4126799488945
$ws.Range("H87").Value = 99958.03107115487
$ws.Range("I87").Value = '202506131600'

$ws.Range("A88").Value = '范式进化投资组合'
$ws.Range("B88").Value = 'HK02899'
$ws.Range("C88").Value = '紫金矿业'
$ws.Range("D88").Value = 1.11
$ws.Range("E88").Value = 55.44578147671601
$ws.Range("F88").Value = 20.1
$ws.Range("G88").Value = 1114.460207681992
$ws.Range("H88").Value = 99958.03107115487
$ws.Range("I88").Value = '202506131600'

$ws.Range("A89").Value = '范式进化投资组合'
$ws.Range("B89").Value = 'HK06881'
$ws.Range("C89").Value = '中国银河'
$ws.Range("D89").Value = 1
$ws.Range("E89").Value = 122.2789606161316
$ws.Range("F89").Value = 8.18
$ws.Range("G89").Value = 1000.241897839956
$ws.Range("H89").Value = 99958.03107115487
$ws.Range("I89").Value = '202506131600'

$ws.Range("A90").Value = '范式进化投资组合'
$ws.Range("B90").Value = '100000'
$ws.Range("C90").Value = '现金'
$ws.Range("D90").Value = 66.04000000000001
$ws.Range("E90").Value = 66015.96525743716
$ws.Range("F90").Value = 1
$ws.Range("G90").Value = 66015.96525743716
$ws.Range("H90").Value = 99958.03107115487
$ws.Range("I90").Value = '202506131600'

